# Update gh-pages to output generated at 456a3b4
# Applies refreshed "想去人数" (interested-count) figures scraped from
# bilibili show pages, plus one venue relocation (新海诚动漫 某某主题餐厅)
# whose address / cover image changed.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- Sheet 展览 (sheet1) ----
$ws1.Range("F3").Value = 14066
$ws1.Range("F4").Value = 14066
$ws1.Range("F5").Value = 14103
$ws1.Range("F7").Value = 1388
$ws1.Range("F8").Value = 5837
$ws1.Range("F9").Value = 975
$ws1.Range("D14").Value = "康候圣街99号 顺丰创新中心"
$ws1.Range("F14").Value = 1524
$ws1.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202407/LHJAvvfX1721095909519.png"
$ws1.Range("F15").Value = 426
$ws1.Range("F16").Value = 2129
$ws1.Range("F17").Value = 1182
$ws1.Range("F18").Value = 1799
$ws1.Range("F19").Value = 910
$ws1.Range("F21").Value = 2256
$ws1.Range("F22").Value = 553
$ws1.Range("F23").Value = 795
$ws1.Range("F24").Value = 3277
$ws1.Range("F26").Value = 305
$ws1.Range("F27").Value = 2342
$ws1.Range("F28").Value = 578
$ws1.Range("F32").Value = 1066
$ws1.Range("F33").Value = 1352
$ws1.Range("F34").Value = 96
$ws1.Range("F36").Value = 4696
$ws1.Range("F37").Value = 4760
$ws1.Range("F40").Value = 666
$ws1.Range("F41").Value = 674
$ws1.Range("F42").Value = 3268
$ws1.Range("F45").Value = 333
$ws1.Range("F46").Value = 89
$ws1.Range("F47").Value = 65
$ws1.Range("F48").Value = 4408
$ws1.Range("F49").Value = 553
$ws1.Range("F50").Value = 276

# ---- Sheet 演出 (sheet2) ----
$ws2.Range("F4").Value = 113

# ---- Sheet 本地生活 (sheet3) ----
$ws3.Range("F2").Value = 7445
$ws3.Range("F4").Value = 696

# ---- Sheet 全部类型 (sheet4) ----
$ws4.Range("F2").Value = 7445
$ws4.Range("F5").Value = 696
$ws4.Range("F6").Value = 14066
$ws4.Range("F7").Value = 14104
$ws4.Range("F9").Value = 1388
$ws4.Range("F10").Value = 5837
$ws4.Range("F11").Value = 975
$ws4.Range("F12").Value = 113
$ws4.Range("D15").Value = "康候圣街99号 顺丰创新中心"
$ws4.Range("F15").Value = 1524
$ws4.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202407/LHJAvvfX1721095909519.png"
$ws4.Range("F16").Value = 426
$ws4.Range("F17").Value = 1182
$ws4.Range("F18").Value = 1799
$ws4.Range("F19").Value = 910
$ws4.Range("F21").Value = 3277
$ws4.Range("F22").Value = 305
$ws4.Range("F23").Value = 2342
$ws4.Range("F24").Value = 578
$ws4.Range("F31").Value = 1066
$ws4.Range("F32").Value = 1352
$ws4.Range("F33").Value = 96
$ws4.Range("F34").Value = 4696
$ws4.Range("F35").Value = 4760
$ws4.Range("F38").Value = 666
$ws4.Range("F39").Value = 674
$ws4.Range("F40").Value = 3268
$ws4.Range("F42").Value = 333
$ws4.Range("F43").Value = 89
$ws4.Range("F45").Value = 65
$ws4.Range("F46").Value = 4408
$ws4.Range("F47").Value = 553
$ws4.Range("F48").Value = 276
